$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New rows to append to the Packages table (columns: A=PackageType, B=State,
# C=Authority, D=ActionType, E=PackageID, F=Status, G=ParentID)
$rows = @(
    @("Waiver", "MD", "1915(c)",      "Amendment", "MD-2260.R00.70", "Submitted",           "MD-2260.R00.00"),
    @("SPA",    "MD", "Medicaid SPA", "",          "MD-25-9554",     "Approved",            ""),
    @("SPA",    "MD", "Medicaid SPA", "",          "MD-25-9555",     "Under Review",        ""),
    @("SPA",    "MD", "Medicaid SPA", "",          "MD-25-9556",     "Disapproved",         ""),
    @("Waiver", "MD", "1915(c)",      "Amendment", "MD-2260.R00.71", "Pending-Approval",    "MD-2260.R00.00"),
    @("SPA",    "MD", "Medicaid SPA", "",          "MD-25-9557",     "Pending-Concurrence", ""),
    @("SPA",    "MD", "CHIP SPA",     "",          "MD-25-9558",     "Submitted",           ""),
    @("SPA",    "MD", "Medicaid SPA", "",          "MD-25-9559",     "RAI Issued",          ""),
    @("SPA",    "MD", "CHIP SPA",     "",          "MD-25-9560",     "Submitted",           ""),
    @("SPA",    "MD", "Medicaid SPA", "",          "MD-25-9561",     "Submitted",           ""),
    @("SPA",    "MD", "Medicaid SPA", "",          "MD-25-9562",     "Submitted",           ""),
    @("Waiver", "MD", "1915(c)",      "Amendment", "MD-2260.R00.72", "Unsubmitted",         "MD-2260.R00.00"),
    @("SPA",    "MD", "Medicaid SPA", "",          "MD-25-9563",     "Under Review",        ""),
    @("Waiver", "MD", "1915(b)",      "Initial",   "MD-2286.R00.00", "Terminated",          ""),
    @("SPA",    "MD", "Medicaid SPA", "",          "MD-25-9564",     "Withdrawn",           "")
)

$startRow = 79
for ($i = 0; $i -lt $rows.Count; $i++) {
    $r = $startRow + $i
    $rowData = $rows[$i]
    for ($c = 0; $c -lt $rowData.Count; $c++) {
        $ws.Cells.Item($r, $c + 1).Value = $rowData[$c]
    }
}
